# New Submission Synced: 2026-02-08 20:32:13
# Appends the new form-response row (row 4) to the "JSS 3D" sheet,
# matching the columns: Timestamp | Full Name | Admission No | AI Score

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("JSS 3D")

$ws.Range("A4").Value = "2026-02-08 20:32:13"
$ws.Range("B4").Value = "Yahya A Yahya"

# "Admission No" is stored as text even though it looks numeric (matches
# the existing rows, e.g. C3 = "38"), so force text via the leading
# apostrophe, then strip the resulting quote-prefix style so the cell
# keeps the plain, unstyled look of its neighbours.
$ws.Range("C4").Value = "'24"
$ws.Range("C4").Style = "Normal"

$ws.Range("D4").Value = 10
